$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 6604472408411690
$ws.Range("D3").Value = 8232260351980812
$ws.Range("D4").Value = 5712944557296251
$ws.Range("D5").Value = 24770837348773640
$ws.Range("D6").Value = 2954685332635526
$ws.Range("D7").Value = 1999376559383782
$ws.Range("D8").Value = 3584075249010192
$ws.Range("D9").Value = 14485860614020270
$ws.Range("D10").Value = 1867272210.335898
$ws.Range("D11").Value = 952044472928671.6
$ws.Range("D12").Value = 1519782252619251
$ws.Range("D13").Value = 882363431162399.4
$ws.Range("D14").Value = 1451515514439668
$ws.Range("D15").Value = 1551416074473771
$ws.Range("D16").Value = 2933666768474808
$ws.Range("D17").Value = 2400430098853130
$ws.Range("D18").Value = 1971728870494816
$ws.Range("D19").Value = 40402161021.95883
